# "Generate Report for handoff"
#
# The localization-status report gets a newly-discovered source file
# (6f4d5cb9-d189-4dc5-b965-ef6ac3bcd797.md) recorded with a "Handoff
# failed" status, inserted as a row ahead of the ".localization-config"
# row on every sheet. The existing source file's generated UUID also
# rolled to a fresh handoff (7f4533f2-41af-4083-aeaf-57f75932345a.md)
# with a new target-file hash/timestamp.

$wb = $excel.ActiveWorkbook

$oldMd  = "9b75be1d-03c0-4178-823f-9885773f5ee1.md"
$newMd  = "7f4533f2-41af-4083-aeaf-57f75932345a.md"
$failMd = "6f4d5cb9-d189-4dc5-b965-ef6ac3bcd797.md"

$oldXlfBase = "9b75be1d-03c0-4178-823f-9885773f5ee1.e55d3901299cfb522614a0c17389f5f2f3d3a799"
$newXlfBase = "7f4533f2-41af-4083-aeaf-57f75932345a.8995df78cfd783cc48509ad836dd4eaf705e7719"

$zhXlf = "$newXlfBase.zh-cn.xlf"
$deXlf = "$newXlfBase.de-de.xlf"

$zhHandoffDt = "2016-01-11 03:23:15"
$deHandoffDt = "2016-01-11 03:23:31"

$epoch = "0001-01-01 00:00:00"

$repoBase    = "https://github.com/OpenLocalizationTest/oltest/blob/2c0a8311d27d69677c73d01793650771dcd2927a"
$zhXlfBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e9c418f4fe9bc089cf62cca570fd66e660177e92/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang"
$deXlfBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d6863cda52a962d7eaf71a7e9602a1878e2b477/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows.Item(3).Insert()

$wsOverview.Range("A3").Value = $failMd
$wsOverview.Range("B3").Value = "Handoff failed"
$wsOverview.Range("C3").Value = "Handoff failed"

$wsOverview.Range("A2").Value = $newMd

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "$repoBase/e2e/$newMd", "", "", $newMd)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "$repoBase/e2e/$failMd", "", "", $failMd)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "$repoBase/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows.Item(3).Insert()

$wsZh.Range("A3").Value = $failMd
$wsZh.Range("B3").Value = "Handoff failed"
$wsZh.Range("D3").Value = $epoch
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = "Ignored"

$wsZh.Range("A2").Value = $newMd
$wsZh.Range("C2").Value = $zhXlf
$wsZh.Range("D2").Value = $zhHandoffDt

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$repoBase/e2e/$newMd", "", "", $newMd)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "$zhXlfBase/$zhXlf", "", "", $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$repoBase/e2e/$failMd", "", "", $failMd)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "$repoBase/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows.Item(3).Insert()

$wsDe.Range("A3").Value = $failMd
$wsDe.Range("B3").Value = "Handoff failed"
$wsDe.Range("D3").Value = $epoch
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = "Ignored"

$wsDe.Range("A2").Value = $newMd
$wsDe.Range("C2").Value = $deXlf
$wsDe.Range("D2").Value = $deHandoffDt

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$repoBase/e2e/$newMd", "", "", $newMd)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "$deXlfBase/$deXlf", "", "", $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$repoBase/e2e/$failMd", "", "", $failMd)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "$repoBase/.localization-config", "", "", ".localization-config")
